# auto-update timeframe-hour MOB --> da comletare: e' solo un test
# Refresh the "Report Activity" metrics (columns B:L) for every category row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, then the new values for columns B..L (in order):
#   B=Unnamed:0  C=Attesa  D=Risposte Eff.  E=Offerte  F=Abb sup.14
#   G=Abb inf.14 H=T.A.    I=Livello di Servizio %  J=Short Call min 10
#   K=Cleared    L=% Cleared
$rows = @(
    @(2,  @(987,  2, 1340, 1567, 47,  14, 4030, 22395, 25, 164, 3369)),  # AMM
    @(3,  @(265,  0, 133,  142,  8,   1,  567,  4846,  2,  0,   0)),     # AMM CT
    @(4,  @(6,    2, 106,  118,  9,   1,  769,  1810,  4,  0,   0)),     # COM
    @(5,  @(1543, 0, 114,  117,  3,   0,  1611, 24975, 0,  0,   0)),     # IPR
    @(6,  @(1813, 0, 251,  251,  2,   0,  1985, 87005, 8,  0,   0)),     # MIG
    @(7,  @(133,  0, 0,    0,    0,   0,  0,    0,     0,  0,   0)),     # MKT
    @(8,  @(1027, 0, 798,  899,  27,  2,  473,  14040, 15, 73,  1164)),  # MOB
    @(9,  @(264,  0, 74,   77,   3,   0,  390,  5000,  1,  0,   0)),     # MOB CT
    @(10, @(576,  5, 490,  631,  130, 12, 6119, 6425,  5,  0,   0)),     # MOB PRE
    @(11, @(1040, 1, 198,  206,  5,   1,  284,  5018,  1,  0,   0)),     # MSK
    @(12, @(1396, 0, 60,   78,   18,  0,  8800, 14576, 0,  0,   0)),     # NOT
    @(13, @(865,  0, 403,  404,  0,   0,  0,    13446, 11, 1,   370)),   # TEC
    @(14, @(266,  0, 146,  151,  4,   1,  267,  5105,  1,  0,   0)),     # TEC CT
    @(15, @(295,  1, 63,   70,   6,   0,  1250, 13770, 0,  0,   0)),     # TST
    @(16, @(187,  0, 1,    1,    0,   0,  0,    0,     0,  0,   0)),     # VIP
    @(17, @(1801, 0, 58,   64,   5,   1,  2335, 24002, 0,  0,   0)),     # WLC
    @(18, @(4079, 3, 430,  981,  48,  3,  1421, 21701, 9,  497, 14496))  # ZERO
)

foreach ($entry in $rows) {
    $r = $entry[0]
    $vals = $entry[1]
    $col = 2  # column B
    foreach ($v in $vals) {
        $ws.Cells.Item($r, $col).Value = $v
        $col = $col + 1
    }
}
